$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2" = 13.16239766666667
    "H2" = 39.487193
    "I2" = 0.5184384809583861
    "J2" = 0.5184384809583861
    "M2" = 11.6425
    "N2" = 34.9275
    "O2" = 0.3149027819311015
    "P2" = 0.3149027819311015
    "Q2" = 153.2432148341667
    "R2" = 1379.1889335075
    "S2" = 0.1632577199139302
    "T2" = 0.1632577199139302
    "G3" = 13.16239766666667
    "H3" = 39.487193
    "I3" = 0.5184384809583861
    "J3" = 0.5184384809583861
    "O3" = 0.5443380502348534
    "P3" = 0.5443380502348534
    "Q3" = 264.8948105920572
    "R3" = 2384.053295328515
    "S3" = 0.282205791891607
    "T3" = 0.282205791891607
    "G4" = 13.16239766666667
    "H4" = 39.487193
    "I4" = 0.5184384809583861
    "J4" = 0.5184384809583861
    "M4" = 5.204109666666667
    "N4" = 15.612329
    "O4" = 0.1407591678340452
    "P4" = 0.1407591678340452
    "Q4" = 68.49856093361079
    "R4" = 616.4870484024971
    "S4" = 0.07297496915284889
    "T4" = 0.07297496915284889
    "I5" = 0.3483468901701054
    "J5" = 0.3483468901701054
    "M5" = 11.6425
    "N5" = 34.9275
    "O5" = 0.3149027819311015
    "P5" = 0.3149027819311015
    "Q5" = 102.96650285
    "R5" = 926.6985256500001
    "S5" = 0.1096954047916141
    "T5" = 0.1096954047916141
    "I6" = 0.3483468901701054
    "J6" = 0.3483468901701054
    "O6" = 0.5443380502348534
    "P6" = 0.5443380502348534
    "S6" = 0.1896184670005698
    "T6" = 0.1896184670005698
    "I7" = 0.3483468901701054
    "J7" = 0.3483468901701054
    "M7" = 5.204109666666667
    "N7" = 15.612329
    "O7" = 0.1407591678340452
    "P7" = 0.1407591678340452
    "Q7" = 46.02524997419334
    "R7" = 414.22724976774
    "S7" = 0.04903301837792156
    "T7" = 0.04903301837792157
    "G8" = 3.382125333333333
    "H8" = 10.146376
    "I8" = 0.1332146288715084
    "J8" = 0.1332146288715084
    "M8" = 11.6425
    "N8" = 34.9275
    "O8" = 0.3149027819311015
    "P8" = 0.3149027819311015
    "Q8" = 39.37639419333333
    "R8" = 354.38754774
    "S8" = 0.04194965722555723
    "T8" = 0.04194965722555724
    "G9" = 3.382125333333333
    "H9" = 10.146376
    "I9" = 0.1332146288715084
    "J9" = 0.1332146288715084
    "O9" = 0.5443380502348534
    "P9" = 0.5443380502348534
    "Q9" = 68.06567255149777
    "R9" = 612.59105296348
    "S9" = 0.07251379134267649
    "T9" = 0.0725137913426765
    "G10" = 3.382125333333333
    "H10" = 10.146376
    "I10" = 0.1332146288715084
    "J10" = 0.1332146288715084
    "M10" = 5.204109666666667
    "N10" = 15.612329
    "O10" = 0.1407591678340452
    "P10" = 0.1407591678340452
    "Q10" = 17.60095114107822
    "R10" = 158.408560269704
    "S10" = 0.01875118030327469
    "T10" = 0.01875118030327469
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
